# Update "paises" workbook:
#  - Arabia Saudita overtakes Japon in total cases -> reorder rows 26-29
#  - Eslovenia overtakes Bosnia y Herzegovina -> reorder rows 76-77
#  - Mali overtakes Tanzania -> reorder rows 122-123
#  - Refresh a handful of other per-country stat rows (37, 115)
#  - Bump the "last updated" timestamp in the footer cell (A1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($row, $values) {
    $cols = @("A","B","C","D","E","F","G","H")
    for ($i = 0; $i -lt $values.Length; $i++) {
        $ws.Range("$($cols[$i])$row").Value = $values[$i]
    }
}

# --- Footer timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 21 de Abril de 2020 a las 14:52"

# --- Rows 26-29: Arabia Saudita moves ahead of Japon -------------------
Set-Row 26 @("Arabia Saudita", 11631, 1147, 1640, 9882, 88, 6, 109)
Set-Row 27 @("Japon", 11135, 0, 1239, 9633, 217, 0, 263)
Set-Row 28 @("Corea del Sur", 10683, 9, 8213, 2233, 55, 1, 237)
Set-Row 29 @("Chile", 10507, 0, 4676, 5692, 377, 0, 139)

# --- Row 37: Dinamarca stats refresh (country stays in place) ----------
$ws.Range("D37").Value = 4700
$ws.Range("E37").Value = 2625
$ws.Range("F37").Value = 72
$ws.Range("G37").Value = 6
$ws.Range("H37").Value = 370

# --- Rows 76-77: Eslovenia moves ahead of Bosnia y Herzegovina ---------
Set-Row 76 @("Eslovenia", 1344, 9, 197, 1070, 25, 0, 77)
Set-Row 77 @("Bosnia y Herzegovina", 1342, 33, 437, 854, 4, 2, 51)

# --- Row 115: Sri Lanka stats refresh (country stays in place) ---------
$ws.Range("B115").Value = 310
$ws.Range("C115").Value = 6
$ws.Range("D115").Value = 102
$ws.Range("E115").Value = 201

# --- Rows 122-123: Mali moves ahead of Tanzania -------------------------
Set-Row 122 @("Mali", 258, 12, 57, 187, 0, 0, 14)
Set-Row 123 @("Tanzania", 254, 0, 11, 233, 4, 0, 10)
